{"js": "// Design notes for branching\n// Insert the word \"completed \" before \"sessions alongside ...\" in the\n// paragraph that introduces the branching timeline figure, so that the\n// sentence reads \"...involving two completed sessions alongside...\".\n\nconst body = context.document.body;\n\n// Locate the unique sentence fragment that needs to change.\nconst searchResults = body.search(\n  \"two sessions alongside the data that would be recorded in the\",\n  { matchCase: false, matchWholeWord: false }\n);\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Could not find target sentence to update.\");\n}\n\nconst hit = searchResults.items[0];\n\n// Narrow down to just the word \"sessions\" inside the matched range so we\n// can insert \"completed \" immediately before it, turning\n// \"...two sessions alongside...\" into \"...two completed sessions alongside...\".\nconst sessionsResults = hit.search(\"sessions alongside\", { matchCase: false });\nsessionsResults.load(\"items\");\nawait context.sync();\n\nif (sessionsResults.items.length === 0) {\n  throw new Error(\"Could not find 'sessions alongside' to update.\");\n}\n\nsessionsResults.items[0].insertText(\"completed \", Word.InsertLocation.before);\nawait context.sync();\n", "ps1": "# Design notes for branching\n# Insert the word \"completed \" before \"sessions alongside ...\" in the\n# paragraph that introduces the branching timeline figure, so that the\n# sentence reads \"...involving two completed sessions alongside...\".\n\n$d = $word.ActiveDocument\n\n$target = $d.Content\n$target.Find.ClearFormatting()\n$target.Find.MatchCase = $false\n$target.Find.Text = \"sessions alongside the data that would be recorded in the\"\n$target.Find.Execute() | Out-Null\n\nif ($target.Find.Found) {\n    $target.Collapse(1)  # wdCollapseStart\n    $target.InsertBefore(\"completed \")\n}\n"}
